$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 855.1875
$ws.Range("J12").Value = 1158.3334
$ws.Range("L12").Value = 1158.3334
$ws.Range("N12").Value = -1498.3334
$ws.Range("H17").Value = 2342.2666
$ws.Range("J17").Value = 2352.7856
$ws.Range("L17").Value = 7058.3568
$ws.Range("N17").Value = -7394.3568
$ws.Range("H33").Value = 904.3684
$ws.Range("J33").Value = 980.5714
$ws.Range("L33").Value = 980.5714
$ws.Range("N33").Value = -1438.5714
$ws.Range("H38").Value = 1459.8889
$ws.Range("I38").Value = 106.5
$ws.Range("K38").Value = 319.5
$ws.Range("M38").Value = 52.5
$ws.Range("H51").Value = 12504305
$ws.Range("I51").Value = 31254124
$ws.Range("J51").Value = 4424.6665
$ws.Range("K51").Value = 31254124
$ws.Range("L51").Value = 4424.6665
$ws.Range("M51").Value = -31253640
$ws.Range("N51").Value = -5392.6665
$ws.Range("H62").Value = 7800.758
$ws.Range("I62").Value = 6954.893
$ws.Range("K62").Value = 6954.893
$ws.Range("M62").Value = -6330.893
$ws.Range("H64").Value = 13817.318
$ws.Range("I64").Value = 16373.8125
$ws.Range("J64").Value = 7000
$ws.Range("K64").Value = 16373.8125
$ws.Range("L64").Value = 7000
$ws.Range("M64").Value = -16125.8125
$ws.Range("N64").Value = -7496
$ws.Range("H65").Value = 7800.758
$ws.Range("I65").Value = 6954.893
$ws.Range("K65").Value = 34774.465
$ws.Range("M65").Value = -31654.465
$ws.Range("H67").Value = 13817.318
$ws.Range("I67").Value = 16373.8125
$ws.Range("J67").Value = 7000
$ws.Range("K67").Value = 16373.8125
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = -15515.8125
$ws.Range("N67").Value = -8716
$ws.Range("H70").Value = 11369
$ws.Range("J70").Value = 8149.75
$ws.Range("L70").Value = 24449.25
$ws.Range("N70").Value = -24989.25
$ws.Range("H73").Value = 11369
$ws.Range("J73").Value = 8149.75
$ws.Range("L73").Value = 24449.25
$ws.Range("N73").Value = -26321.25
$ws.Range("H76").Value = 4663.0835
$ws.Range("I76").Value = 3870.25
$ws.Range("K76").Value = 3870.25
$ws.Range("M76").Value = -3555.25
$ws.Range("H79").Value = 4663.0835
$ws.Range("I79").Value = 3870.25
$ws.Range("K79").Value = 3870.25
$ws.Range("M79").Value = -2778.25
$ws.Range("H94").Value = 1878.125
$ws.Range("I94").Value = 1878.125
$ws.Range("K94").Value = 1878.125
$ws.Range("M94").Value = -1427.125
$ws.Range("H99").Value = 511.66666
$ws.Range("I99").Value = 321.33334
$ws.Range("K99").Value = 964.0000200000001
$ws.Range("M99").Value = 533.9999799999999
$ws.Range("H100").Value = 52961
$ws.Range("I100").Value = 52961
$ws.Range("K100").Value = 52961
$ws.Range("M100").Value = -52420
$ws.Range("H112").Value = 6695.4224
$ws.Range("J112").Value = 6813.5
$ws.Range("L112").Value = 20440.5
$ws.Range("N112").Value = -22656.5
$ws.Range("H116").Value = 5009750
$ws.Range("I116").Value = 6677999.5
$ws.Range("K116").Value = 6677999.5
$ws.Range("M116").Value = -6674557.5
$ws.Range("H132").Value = 2762794
$ws.Range("I132").Value = 2816955
$ws.Range("J132").Value = 584
$ws.Range("K132").Value = 8450865
$ws.Range("L132").Value = 1752
$ws.Range("M132").Value = -8448335
$ws.Range("N132").Value = -6812
$ws.Range("H138").Value = 1839.2941
$ws.Range("J138").Value = 3481.5
$ws.Range("L138").Value = 10444.5
$ws.Range("N138").Value = -20724.5
$ws.Range("H141").Value = 2307.8857
$ws.Range("I141").Value = 1910.1034
$ws.Range("J141").Value = 4230.5
$ws.Range("K141").Value = 5730.3102
$ws.Range("L141").Value = 12691.5
$ws.Range("M141").Value = -550.3101999999999
$ws.Range("N141").Value = -23051.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2309.5
$ws.Range("I2").Value = 1940
$ws.Range("J2").Value = 2956.125
$ws.Range("K2").Value = 1940
$ws.Range("L2").Value = 2956.125
$ws.Range("M2").Value = -1827
$ws.Range("N2").Value = -3182.125
$ws.Range("H32").Value = 24402.096
$ws.Range("I32").Value = 26867.37
$ws.Range("K32").Value = 26867.37
$ws.Range("M32").Value = -26580.37
$ws.Range("H36").Value = 19900
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 10000
$ws.Range("K42").Value = 10000
$ws.Range("M42").Value = -9514
$ws.Range("H55").Value = 27997.334
$ws.Range("J55").Value = 27997.334
$ws.Range("L55").Value = 27997.334
$ws.Range("N55").Value = -28627.334
$ws.Range("H61").Value = 5231.1377
$ws.Range("I61").Value = 1088.2858
$ws.Range("K61").Value = 1088.2858
$ws.Range("M61").Value = -876.2858000000001
$ws.Range("H74").Value = 198012.03
$ws.Range("I74").Value = 241215.8
$ws.Range("K74").Value = 241215.8
$ws.Range("M74").Value = -240341.8
$ws.Range("H77").Value = 198012.03
$ws.Range("I77").Value = 241215.8
$ws.Range("K77").Value = 1206079
$ws.Range("M77").Value = -1201711
$ws.Range("H102").Value = 7101.6665
$ws.Range("I102").Value = 5522
$ws.Range("K102").Value = 5522
$ws.Range("M102").Value = -3900
$ws.Range("H110").Value = 1992.6364
$ws.Range("I110").Value = 2047.6666
$ws.Range("J110").Value = 1745
$ws.Range("K110").Value = 2047.6666
$ws.Range("L110").Value = 1745
$ws.Range("M110").Value = -2.666600000000017
$ws.Range("N110").Value = -5835
$ws.Range("H116").Value = 2309.5
$ws.Range("I116").Value = 1940
$ws.Range("J116").Value = 2956.125
$ws.Range("K116").Value = 1940
$ws.Range("L116").Value = 2956.125
$ws.Range("M116").Value = 354
$ws.Range("N116").Value = -7544.125
$ws.Range("H132").Value = 1292.2273
$ws.Range("I132").Value = 970.3699
$ws.Range("J132").Value = 2858.6
$ws.Range("K132").Value = 2911.1097
$ws.Range("L132").Value = 8575.799999999999
$ws.Range("M132").Value = -381.1097
$ws.Range("N132").Value = -13635.8
$ws.Range("H136").Value = 5231.1377
$ws.Range("I136").Value = 1088.2858
$ws.Range("K136").Value = 3264.8574
$ws.Range("M136").Value = -714.8574000000003
$ws.Range("M36").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2309.5
$ws.Range("I3").Value = 1940
$ws.Range("J3").Value = 2956.125
$ws.Range("K3").Value = 1940
$ws.Range("L3").Value = 2956.125
$ws.Range("M3").Value = -1826
$ws.Range("N3").Value = -3184.125
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("H74").Value = 73963.336
$ws.Range("J74").Value = 73963.336
$ws.Range("L74").Value = 73963.336
$ws.Range("N74").Value = -75835.336
$ws.Range("H77").Value = 73963.336
$ws.Range("J77").Value = 73963.336
$ws.Range("L77").Value = 221890.008
$ws.Range("N77").Value = -231250.008
$ws.Range("H86").Value = 1244.24
$ws.Range("I86").Value = 1222.125
$ws.Range("J86").Value = 1283.5555
$ws.Range("K86").Value = 1222.125
$ws.Range("L86").Value = 1283.5555
$ws.Range("M86").Value = -99.125
$ws.Range("N86").Value = -3529.5555
$ws.Range("H89").Value = 1244.24
$ws.Range("I89").Value = 1222.125
$ws.Range("J89").Value = 1283.5555
$ws.Range("K89").Value = 6110.625
$ws.Range("L89").Value = 6417.7775
$ws.Range("M89").Value = -494.625
$ws.Range("N89").Value = -17649.7775
$ws.Range("H94").Value = 2216.0557
$ws.Range("I94").Value = 1138.1111
$ws.Range("J94").Value = 3294
$ws.Range("K94").Value = 1138.1111
$ws.Range("L94").Value = 3294
$ws.Range("M94").Value = -687.1111000000001
$ws.Range("N94").Value = -4196
$ws.Range("H99").Value = 2577.6667
$ws.Range("I99").Value = 2342.7144
$ws.Range("K99").Value = 2342.7144
$ws.Range("M99").Value = -844.7143999999998
$ws.Range("H105").Value = 3526.55
$ws.Range("I105").Value = 3666.853
$ws.Range("J105").Value = 2731.5
$ws.Range("K105").Value = 3666.853
$ws.Range("L105").Value = 2731.5
$ws.Range("M105").Value = -1919.853
$ws.Range("N105").Value = -6225.5
$ws.Range("H107").Value = 17201.06
$ws.Range("I107").Value = 24072
$ws.Range("K107").Value = 24072
$ws.Range("M107").Value = -22152
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H111").Value = 95000
$ws.Range("J111").Value = 95000
$ws.Range("L111").Value = 95000
$ws.Range("N111").Value = -103180
$ws.Range("H134").Value = 2211.3767
$ws.Range("I134").Value = 1293.8551
$ws.Range("K134").Value = 3881.5653
$ws.Range("M134").Value = -1346.5653
$ws.Range("H137").Value = 118328.57
$ws.Range("J137").Value = 118328.57
$ws.Range("L137").Value = 118328.57
$ws.Range("N137").Value = -128528.57
$ws.Range("H139").Value = 114479.82
$ws.Range("J139").Value = 114479.82
$ws.Range("L139").Value = 114479.82
$ws.Range("N139").Value = -124759.82
$ws.Range("H140").Value = 102493.6
$ws.Range("J140").Value = 102493.6
$ws.Range("L140").Value = 102493.6
$ws.Range("N140").Value = -112853.6
$ws.Range("N61").ClearContents()
$ws.Range("N109").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 165237500
$ws.Range("I12").Value = 225156260
$ws.Range("J12").Value = 45400000
$ws.Range("K12").Value = 225156260
$ws.Range("L12").Value = 45400000
$ws.Range("M12").Value = -225156090
$ws.Range("N12").Value = -45400340
$ws.Range("H16").Value = 1061.25
$ws.Range("I16").Value = 966.875
$ws.Range("K16").Value = 966.875
$ws.Range("M16").Value = -679.875
$ws.Range("H19").Value = 3649.8
$ws.Range("I19").Value = 928.2857
$ws.Range("K19").Value = 928.2857
$ws.Range("M19").Value = -758.2857
$ws.Range("H24").Value = 3649.8
$ws.Range("I24").Value = 928.2857
$ws.Range("K24").Value = 928.2857
$ws.Range("M24").Value = -758.2857
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("H31").Value = 5265715.5
$ws.Range("I31").Value = 7694211.5
$ws.Range("J31").Value = 3975
$ws.Range("K31").Value = 7694211.5
$ws.Range("L31").Value = 3975
$ws.Range("M31").Value = -7693916.5
$ws.Range("N31").Value = -4565
$ws.Range("H34").Value = 5265715.5
$ws.Range("I34").Value = 7694211.5
$ws.Range("J34").Value = 3975
$ws.Range("K34").Value = 7694211.5
$ws.Range("L34").Value = 3975
$ws.Range("M34").Value = -7694009.5
$ws.Range("N34").Value = -4379
$ws.Range("H58").Value = 810
$ws.Range("I58").Value = 833.6286
$ws.Range("K58").Value = 833.6286
$ws.Range("M58").Value = -630.6286
$ws.Range("H62").Value = 5696.4287
$ws.Range("I62").Value = 3186
$ws.Range("J62").Value = 7978.636
$ws.Range("K62").Value = 3186
$ws.Range("L62").Value = 7978.636
$ws.Range("M62").Value = -2562
$ws.Range("N62").Value = -9226.636
$ws.Range("H65").Value = 5696.4287
$ws.Range("I65").Value = 3186
$ws.Range("J65").Value = 7978.636
$ws.Range("K65").Value = 15930
$ws.Range("L65").Value = 39893.18
$ws.Range("M65").Value = -12810
$ws.Range("N65").Value = -46133.18
$ws.Range("H76").Value = 111118490
$ws.Range("I76").Value = 111118490
$ws.Range("K76").Value = 111118490
$ws.Range("M76").Value = -111118175
$ws.Range("H79").Value = 111118490
$ws.Range("I79").Value = 111118490
$ws.Range("K79").Value = 111118490
$ws.Range("M79").Value = -111117398
$ws.Range("H86").Value = 29189.4
$ws.Range("I86").Value = 7950
$ws.Range("J86").Value = 34499.25
$ws.Range("K86").Value = 7950
$ws.Range("L86").Value = 34499.25
$ws.Range("M86").Value = -6827
$ws.Range("N86").Value = -36745.25
$ws.Range("H89").Value = 29189.4
$ws.Range("I89").Value = 7950
$ws.Range("J89").Value = 34499.25
$ws.Range("K89").Value = 39750
$ws.Range("L89").Value = 172496.25
$ws.Range("M89").Value = -34134
$ws.Range("N89").Value = -183728.25
$ws.Range("H99").Value = 5058.6313
$ws.Range("I99").Value = 4264.6665
$ws.Range("K99").Value = 4264.6665
$ws.Range("M99").Value = -2766.6665
$ws.Range("H105").Value = 2304
$ws.Range("I105").Value = 1629.75
$ws.Range("K105").Value = 1629.75
$ws.Range("M105").Value = 117.25
$ws.Range("H107").Value = 3197
$ws.Range("I107").Value = 294.5
$ws.Range("J107").Value = 5519
$ws.Range("K107").Value = 294.5
$ws.Range("L107").Value = 5519
$ws.Range("M107").Value = 1625.5
$ws.Range("N107").Value = -9359
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("H113").Value = 1061.25
$ws.Range("I113").Value = 966.875
$ws.Range("K113").Value = 966.875
$ws.Range("M113").Value = 1203.125
$ws.Range("H126").Value = 5058.6313
$ws.Range("I126").Value = 4264.6665
$ws.Range("K126").Value = 12793.9995
$ws.Range("M126").Value = -10323.9995
$ws.Range("H132").Value = 20902.129
$ws.Range("I132").Value = 26147.834
$ws.Range("J132").Value = 2916.8572
$ws.Range("K132").Value = 78443.50199999999
$ws.Range("L132").Value = 8750.571599999999
$ws.Range("M132").Value = -75913.50199999999
$ws.Range("N132").Value = -13810.5716
$ws.Range("H134").Value = 1817.8846
$ws.Range("I134").Value = 1299.6666
$ws.Range("K134").Value = 3898.9998
$ws.Range("M134").Value = -1363.9998
$ws.Range("H136").Value = 810
$ws.Range("I136").Value = 833.6286
$ws.Range("K136").Value = 2500.8858
$ws.Range("M136").Value = 49.11419999999998
$ws.Range("N28").ClearContents()
$ws.Range("N112").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 452.58823
$ws.Range("I2").Value = 344.46155
$ws.Range("J2").Value = 804
$ws.Range("K2").Value = 2066.7693
$ws.Range("L2").Value = 4824
$ws.Range("M2").Value = -1953.7693
$ws.Range("N2").Value = -5050
$ws.Range("H5").Value = 2494
$ws.Range("I5").Value = 2327.1667
$ws.Range("J5").Value = 2994.5
$ws.Range("K5").Value = 6981.500100000001
$ws.Range("L5").Value = 8983.5
$ws.Range("M5").Value = -6869.500100000001
$ws.Range("N5").Value = -9207.5
$ws.Range("H22").Value = 4452.878
$ws.Range("J22").Value = 4526.7
$ws.Range("L22").Value = 13580.1
$ws.Range("N22").Value = -13918.1
$ws.Range("H23").Value = 302.5
$ws.Range("J23").Value = 365.75
$ws.Range("L23").Value = 1097.25
$ws.Range("N23").Value = -1567.25
$ws.Range("H27").Value = 4452.878
$ws.Range("J27").Value = 4526.7
$ws.Range("L27").Value = 13580.1
$ws.Range("N27").Value = -13784.1
$ws.Range("H81").Value = 4631.778
$ws.Range("I81").Value = 1686
$ws.Range("K81").Value = 5058
$ws.Range("M81").Value = -3935
$ws.Range("H84").Value = 4631.778
$ws.Range("I84").Value = 1686
$ws.Range("K84").Value = 15174
$ws.Range("M84").Value = -9558
$ws.Range("H110").Value = 8881.75
$ws.Range("I110").Value = 5175.6665
$ws.Range("J110").Value = 20000
$ws.Range("K110").Value = 15526.9995
$ws.Range("L110").Value = 60000
$ws.Range("M110").Value = -11436.9995
$ws.Range("N110").Value = -68180
$ws.Range("H113").Value = 2555
$ws.Range("I113").Value = 882
$ws.Range("J113").Value = 2973.25
$ws.Range("K113").Value = 2646
$ws.Range("L113").Value = 8919.75
$ws.Range("M113").Value = -476
$ws.Range("N113").Value = -13259.75
$ws.Range("H135").Value = 2494
$ws.Range("I135").Value = 2327.1667
$ws.Range("J135").Value = 2994.5
$ws.Range("K135").Value = 20944.5003
$ws.Range("L135").Value = 26950.5
$ws.Range("M135").Value = -18409.5003
$ws.Range("N135").Value = -32020.5
$ws.Range("H136").Value = 2654.3635
$ws.Range("I136").Value = 2333.111
$ws.Range("K136").Value = 6999.333
$ws.Range("M136").Value = -1899.333

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1255.4166
$ws.Range("I2").Value = 1915.2
$ws.Range("J2").Value = 784.1429000000001
$ws.Range("K2").Value = 1915.2
$ws.Range("L2").Value = 784.1429000000001
$ws.Range("M2").Value = -1802.2
$ws.Range("N2").Value = -1010.1429
$ws.Range("H26").Value = 16000
$ws.Range("I26").Value = 9000
$ws.Range("J26").Value = 19500
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 19500
$ws.Range("M26").Value = -8720
$ws.Range("N26").Value = -20060
$ws.Range("H45").Value = 19270.75
$ws.Range("I45").Value = 16333.333
$ws.Range("J45").Value = 21033.2
$ws.Range("K45").Value = 16333.333
$ws.Range("L45").Value = 21033.2
$ws.Range("M45").Value = -15774.333
$ws.Range("N45").Value = -22151.2
$ws.Range("H50").Value = 16000
$ws.Range("I50").Value = 9000
$ws.Range("J50").Value = 19500
$ws.Range("K50").Value = 9000
$ws.Range("L50").Value = 19500
$ws.Range("M50").Value = -8502
$ws.Range("N50").Value = -20496
$ws.Range("H70").Value = 5477.375
$ws.Range("I70").Value = 5450.1665
$ws.Range("K70").Value = 5450.1665
$ws.Range("M70").Value = -5180.1665
$ws.Range("H73").Value = 5477.375
$ws.Range("I73").Value = 5450.1665
$ws.Range("K73").Value = 5450.1665
$ws.Range("M73").Value = -4514.1665
$ws.Range("H80").Value = 3965
$ws.Range("I80").Value = 2127
$ws.Range("J80").Value = 9479
$ws.Range("K80").Value = 2127
$ws.Range("L80").Value = 9479
$ws.Range("M80").Value = -1129
$ws.Range("N80").Value = -11475
$ws.Range("H83").Value = 3965
$ws.Range("I83").Value = 2127
$ws.Range("J83").Value = 9479
$ws.Range("K83").Value = 10635
$ws.Range("L83").Value = 47395
$ws.Range("M83").Value = -5643
$ws.Range("N83").Value = -57379
$ws.Range("H97").Value = 2228.25
$ws.Range("I97").Value = 1610.5
$ws.Range("J97").Value = 3463.75
$ws.Range("K97").Value = 1610.5
$ws.Range("L97").Value = 3463.75
$ws.Range("M97").Value = -1114.5
$ws.Range("N97").Value = -4455.75
$ws.Range("H102").Value = 19194.55
$ws.Range("I102").Value = 20445.482
$ws.Range("K102").Value = 20445.482
$ws.Range("M102").Value = -18823.482
$ws.Range("H110").Value = 125000
$ws.Range("J110").Value = 125000
$ws.Range("L110").Value = 125000
$ws.Range("N110").Value = -133180
$ws.Range("H122").Value = 2371.5356
$ws.Range("I122").Value = 2219.077
$ws.Range("J122").Value = 4353.5
$ws.Range("K122").Value = 6657.231000000001
$ws.Range("L122").Value = 13060.5
$ws.Range("M122").Value = -4207.231000000001
$ws.Range("N122").Value = -17960.5
$ws.Range("H132").Value = 2835.7144
$ws.Range("I132").Value = 2426.0908
$ws.Range("K132").Value = 7278.2724
$ws.Range("M132").Value = -4748.2724
$ws.Range("H136").Value = 15558.091
$ws.Range("J136").Value = 15558.091
$ws.Range("L136").Value = 46674.273
$ws.Range("N136").Value = -51774.273

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2299.6
$ws.Range("I31").Value = 2750
$ws.Range("K31").Value = 2750
$ws.Range("M31").Value = -2502
$ws.Range("H46").Value = 4803.5654
$ws.Range("I46").Value = 2214.5715
$ws.Range("J46").Value = 5936.25
$ws.Range("K46").Value = 2214.5715
$ws.Range("L46").Value = 5936.25
$ws.Range("M46").Value = -2026.5715
$ws.Range("N46").Value = -6312.25
$ws.Range("H55").Value = 1190.2084
$ws.Range("I55").Value = 94.09090999999999
$ws.Range("K55").Value = 94.09090999999999
$ws.Range("M55").Value = 78.90909000000001
$ws.Range("H93").Value = 1823.5526
$ws.Range("J93").Value = 1653.4615
$ws.Range("L93").Value = 1653.4615
$ws.Range("N93").Value = -4149.461499999999
$ws.Range("H100").Value = 2442.5715
$ws.Range("I100").Value = 2024.5
$ws.Range("K100").Value = 2024.5
$ws.Range("M100").Value = -1483.5
$ws.Range("H122").Value = 3161.9644
$ws.Range("J122").Value = 3064.6667
$ws.Range("L122").Value = 9194.000100000001
$ws.Range("N122").Value = -14094.0001
$ws.Range("H132").Value = 2891.8845
$ws.Range("I132").Value = 2716.3333
$ws.Range("K132").Value = 8148.999899999999
$ws.Range("M132").Value = -5618.999899999999
$ws.Range("H136").Value = 2877.9492
$ws.Range("I136").Value = 2492.7886
$ws.Range("K136").Value = 7478.3658
$ws.Range("M136").Value = -4928.3658

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 14383
$ws.Range("I28").Value = 9816.25
$ws.Range("K28").Value = 9816.25
$ws.Range("M28").Value = -9468.25
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 15000
$ws.Range("K39").Value = 15000
$ws.Range("M39").Value = -14587
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H126").Value = 627923.1
$ws.Range("I126").Value = 2275.2
$ws.Range("K126").Value = 6825.599999999999
$ws.Range("M126").Value = -4355.599999999999
$ws.Range("H132").Value = 4645681.5
$ws.Range("I132").Value = 7166239.5
$ws.Range("J132").Value = 2548.0527
$ws.Range("K132").Value = 21498718.5
$ws.Range("L132").Value = 7644.158100000001
$ws.Range("M132").Value = -21496188.5
$ws.Range("N132").Value = -12704.1581
$ws.Range("H136").Value = 9418.229499999999
$ws.Range("I136").Value = 10481.794
$ws.Range("J136").Value = 3326.9092
$ws.Range("K136").Value = 31445.382
$ws.Range("L136").Value = 9980.7276
$ws.Range("M136").Value = -28895.382
$ws.Range("N136").Value = -15080.7276
$ws.Range("N110").ClearContents()
